# Remove the long list of childless "PUMP:..." tag paragraphs that follow
# the introductory bold paragraph, leaving just the title and intro text.

$d = $word.ActiveDocument

# Locate the first paragraph whose text begins with "PUMP" (the start of
# the tag list) and the last paragraph in the document (the end of the
# tag list, right before the final section properties).
$startIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "PUMP*") {
        $startIndex = $i
        break
    }
}

if ($startIndex -gt 0) {
    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($d.Paragraphs.Count)

    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
